$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: mark the new "checkout--> currentuser" entry as complete ---
$ws.Range("C22").Value = "complete"

# --- Row 27: flesh out the "todo" status with the real progress note, ---
# --- and recolor it the same orange used by the "inprogress" row (C19) ---
$ws.Range("C27").Value = "in progress need to write function to send isadmin details after clicking on cart button  "
$ws.Range("C27").Interior.Color = $ws.Range("C19").Interior.Color

# --- Column C needs to be much wider now that it holds a long note ---
$ws.Columns.Item(3).ColumnWidth = 71.83

# --- Move the saved cursor/selection to C31 ---
[void]$ws.Range("C31").Select()
